$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-25 Sunday" "2024-02-26 Monday"
Replace-Text "19×56=" "80×46="
Replace-Text "89×99=" "22×72="
Replace-Text "74×76=" "67×87="
Replace-Text "63×29=" "44×38="
Replace-Text "15×29=" "53×73="
Replace-Text "41×53=" "26×40="
Replace-Text "23×88=" "40×89="
Replace-Text "24×50=" "94×38="
Replace-Text "12×81=" "80×73="
Replace-Text "98×25=" "11×65="
Replace-Text "99×27=" "80×38="
Replace-Text "20×87=" "82×59="
Replace-Text "86×56=" "41×59="
Replace-Text "57×46=" "65×27="
Replace-Text "34×23=" "73×61="
Replace-Text "68×28=" "98×60="
Replace-Text "95×82=" "71×18="
Replace-Text "24×21=" "84×55="
Replace-Text "43×53=" "32×89="
Replace-Text "63×25=" "36×23="
Replace-Text "89×33=" "26×34="
Replace-Text "36×86=" "47×43="
Replace-Text "52×69=" "93×67="
Replace-Text "66×71=" "58×37="
Replace-Text "55×49=" "99×14="
